# Applies the "last minute updates" edit to the SUBPART 5337.5 cover
# paragraph:
#   1. Add a thin paragraph border (5-twip spacing on all four sides).
#   2. Increase the paragraph's left indent from 120 -> 225 twips
#      (6pt -> 11.25pt).
#   3. Rename the placeholder ID from
#      **ID__AFFARS_5337_topic_9__ID** to
#      **ID__AFFARS_SUBPART_5337_5__ID**, and drop the stray trailing
#      space run that followed it.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# --- Remove the trailing " " run that follows the placeholder text ---
# The paragraph's text is "**ID__AFFARS_5337_topic_9__ID** \r"; strip the
# paragraph mark, then trim any trailing whitespace and delete it from
# the document so the stray space-only run disappears entirely.
$withMark = $p.Range.Text
$visible = $withMark.Substring(0, $withMark.Length - 1)
$trimmed = $visible.TrimEnd()
if ($trimmed.Length -lt $visible.Length) {
    $delStart = $p.Range.Start + $trimmed.Length
    $delEnd = $p.Range.Start + $visible.Length
    $d.Range($delStart, $delEnd).Delete()
}

# --- Paragraph formatting: left indent + four-sided border ---
$p.Format.LeftIndent = 11.25

$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5

# --- Rename the placeholder ID text ---
$d.Content.Find.Execute("**ID__AFFARS_5337_topic_9__ID**", $true, $false,
                         $false, $false, $false, $true, 1, $false,
                         "**ID__AFFARS_SUBPART_5337_5__ID**", 2)
